$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 1: add P1 and Q1 (copy O1 header formatting: bold, centered, bordered)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Build a 2D array for B2:Q25 (24 rows x 16 cols: B..Q)
$arr = New-Object 'object[,]' 24,16
$arr[0,0] = 2.849566727709032
$arr[0,1] = 0.9823292823700172
$arr[0,2] = 0.1311576868928483
$arr[0,3] = 0.1516935162507593
$arr[0,4] = 1.268519604685252
$arr[0,5] = 0.9322881981361206
$arr[0,6] = 0
$arr[0,7] = 0.00334714836598371
$arr[0,8] = 0.5876796403362476
$arr[0,9] = 0.4129079284310819
$arr[0,10] = 0.2030316111607391
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0.6985762884866844
$arr[0,15] = 0
$arr[1,0] = 2.477048397991041
$arr[1,1] = 0.8506322728272835
$arr[1,2] = 0.1170531764141458
$arr[1,3] = 0.1352279115194044
$arr[1,4] = 1.165576689919178
$arr[1,5] = 0.8498046374439951
$arr[1,6] = 0
$arr[1,7] = 0.001581032242134572
$arr[1,8] = 0.5558042347527561
$arr[1,9] = 0.3992747321024765
$arr[1,10] = 0.1808443596350671
$arr[1,11] = 0
$arr[1,12] = 0
$arr[1,13] = 0
$arr[1,14] = 0.725143280845014
$arr[1,15] = 0
$arr[2,0] = 2.248328934077108
$arr[2,1] = 0.770679244283258
$arr[2,2] = 0.1085510218632919
$arr[2,3] = 0.1253130076839248
$arr[2,4] = 1.104110424046553
$arr[2,5] = 0.8004869708520062
$arr[2,6] = 0
$arr[2,7] = 0.0008925725355912739
$arr[2,8] = 0.5370205154761578
$arr[2,9] = 0.3914600504199761
$arr[2,10] = 0.1674964651472166
$arr[2,11] = 0
$arr[2,12] = 0
$arr[2,13] = 0
$arr[2,14] = 0.7423461885277582
$arr[2,15] = 0
$arr[3,0] = 2.15358302133393
$arr[3,1] = 0.7390555300712549
$arr[3,2] = 0.1052360833226942
$arr[3,3] = 0.1213524011362779
$arr[3,4] = 1.07851564885317
$arr[3,5] = 0.7797409117766136
$arr[3,6] = 0
$arr[3,7] = 0.0007656690411055322
$arr[3,8] = 0.5290694869257777
$arr[3,9] = 0.3878448440392823
$arr[3,10] = 0.1620955134847364
$arr[3,11] = 0
$arr[3,12] = 0
$arr[3,13] = 0
$arr[3,14] = 0.7497902656098319
$arr[3,15] = 0
$arr[4,0] = 2.135997524591517
$arr[4,1] = 0.7347430057286886
$arr[4,2] = 0.104825868372771
$arr[4,3] = 0.120739972254448
$arr[4,4] = 1.073128756975251
$arr[4,5] = 0.7751494136479948
$arr[4,6] = 0
$arr[4,7] = 0.0008322642982978579
$arr[4,8] = 0.5271721806146843
$arr[4,9] = 0.3865671900540981
$arr[4,10] = 0.1611723804248726
$arr[4,11] = 0
$arr[4,12] = 0
$arr[4,13] = 0
$arr[4,14] = 0.7513075027314251
$arr[4,15] = 0
$arr[5,0] = 2.241975606207632
$arr[5,1] = 0.7727898527983541
$arr[5,2] = 0.1088835910456396
$arr[5,3] = 0.1253758450538562
$arr[5,4] = 1.100588945889385
$arr[5,5] = 0.7970196904892646
$arr[5,6] = 0
$arr[5,7] = 0.001097937434108687
$arr[5,8] = 0.5353050443289504
$arr[5,9] = 0.3895376345200461
$arr[5,10] = 0.1673412438900144
$arr[5,11] = 0
$arr[5,12] = 0
$arr[5,13] = 0
$arr[5,14] = 0.7431945305205652
$arr[5,15] = 0
$arr[6,0] = 2.714292027648469
$arr[6,1] = 0.9400994361407697
$arr[6,2] = 0.1267615811609915
$arr[6,3] = 0.1461273290269389
$arr[6,4] = 1.228414399388129
$arr[6,5] = 0.8993093399895002
$arr[6,6] = 0
$arr[6,7] = 0.002887993508197795
$arr[6,8] = 0.5743809256495069
$arr[6,9] = 0.4055997508901328
$arr[6,10] = 0.1952098595875782
$arr[6,11] = 0
$arr[6,12] = 0
$arr[6,13] = 0
$arr[6,14] = 0.7086128420122932
$arr[6,15] = 0
$arr[7,0] = 3.651379022650474
$arr[7,1] = 1.272170075324766
$arr[7,2] = 0.1625255884253818
$arr[7,3] = 0.18835761639086
$arr[7,4] = 1.500386588613992
$arr[7,5] = 1.117681895478441
$arr[7,6] = 0
$arr[7,7] = 0.00959720247583462
$arr[7,8] = 0.6608377510582102
$arr[7,9] = 0.445362317920484
$arr[7,10] = 0.2525017629959905
$arr[7,11] = 0
$arr[7,12] = 0
$arr[7,13] = 0
$arr[7,14] = 0.6467461160325101
$arr[7,15] = 0
$arr[8,0] = 4.340646402120342
$arr[8,1] = 1.524705295583999
$arr[8,2] = 0.1900322834441113
$arr[8,3] = 0.2208044478093072
$arr[8,4] = 1.710761847299324
$arr[8,5] = 1.286103425980173
$arr[8,6] = 0
$arr[8,7] = 0.01726339619682538
$arr[8,8] = 0.7288142823014709
$arr[8,9] = 0.4773600816607129
$arr[8,10] = 0.2965353833036062
$arr[8,11] = 0
$arr[8,12] = 0
$arr[8,13] = 0
$arr[8,14] = 0.6068075729958409
$arr[8,15] = 0
$arr[9,0] = 4.646768807846001
$arr[9,1] = 1.645959040102298
$arr[9,2] = 0.2034499019487299
$arr[9,3] = 0.2361240395691908
$arr[9,4] = 1.804421191400408
$arr[9,5] = 1.35999963123092
$arr[9,6] = 0
$arr[9,7] = 0.02172320840545172
$arr[9,8] = 0.758455394465102
$arr[9,9] = 0.4898039783473038
$arr[9,10] = 0.3169660353883756
$arr[9,11] = 0
$arr[9,12] = 0
$arr[9,13] = 0
$arr[9,14] = 0.5914609744599986
$arr[9,15] = 0
$arr[10,0] = 4.767572909546971
$arr[10,1] = 1.689948941113585
$arr[10,2] = 0.2082402716379619
$arr[10,3] = 0.2418836774557818
$arr[10,4] = 1.843246826157369
$arr[10,5] = 1.39124966036934
$arr[10,6] = 0
$arr[10,7] = 0.02337793763511264
$arr[10,8] = 0.7713240955293656
$arr[10,9] = 0.4963189619394228
$arr[10,10] = 0.3248680545683698
$arr[10,11] = 0
$arr[10,12] = 0
$arr[10,13] = 0
$arr[10,14] = 0.5850125238728694
$arr[10,15] = 0
$arr[11,0] = 4.742481165418098
$arr[11,1] = 1.680000321034299
$arr[11,2] = 0.2071391307102459
$arr[11,3] = 0.2406201271347328
$arr[11,4] = 1.83543051106605
$arr[11,5] = 1.385071822867815
$arr[11,6] = 0
$arr[11,7] = 0.02298479860685632
$arr[11,8] = 0.7688285316000645
$arr[11,9] = 0.4952395920395389
$arr[11,10] = 0.3231776659992818
$arr[11,11] = 0
$arr[11,12] = 0
$arr[11,13] = 0
$arr[11,14] = 0.586222847282059
$arr[11,15] = 0
$arr[12,0] = 4.657099830240668
$arr[12,1] = 1.649376235526006
$arr[12,2] = 0.2038144672109183
$arr[12,3] = 0.2365880067463237
$arr[12,4] = 1.807846363956742
$arr[12,5] = 1.362804598150575
$arr[12,6] = 0
$arr[12,7] = 0.02184322990048315
$arr[12,8] = 0.7596310131675068
$arr[12,9] = 0.4904772334736762
$arr[12,10] = 0.3176208976757948
$arr[12,11] = 0
$arr[12,12] = 0
$arr[12,13] = 0
$arr[12,14] = 0.5908547048644266
$arr[12,15] = 0
$arr[13,0] = 4.602994620898414
$arr[13,1] = 1.631566657257395
$arr[13,2] = 0.2019165845175337
$arr[13,3] = 0.2341662865139611
$arr[13,4] = 1.789900902308744
$arr[13,5] = 1.348098021189202
$arr[13,6] = 0
$arr[13,7] = 0.02122608659687586
$arr[13,8] = 0.7534643207593206
$arr[13,9] = 0.4869305570910925
$arr[13,10] = 0.3141985804326879
$arr[13,11] = 0
$arr[13,12] = 0
$arr[13,13] = 0
$arr[13,14] = 0.5940491813051878
$arr[13,15] = 0
$arr[14,0] = 4.305217214858942
$arr[14,1] = 1.524409565678354
$arr[14,2] = 0.1902774918015808
$arr[14,3] = 0.2201492756062109
$arr[14,4] = 1.695284484827553
$arr[14,5] = 1.271866350171507
$arr[14,6] = 0
$arr[14,7] = 0.01750664343148856
$arr[14,8] = 0.7221686285746785
$arr[14,9] = 0.4710828115464878
$arr[14,10] = 0.2949566301099367
$arr[14,11] = 0
$arr[14,12] = 0
$arr[14,13] = 0
$arr[14,14] = 0.6105849247288972
$arr[14,15] = 0
$arr[15,0] = 4.123004771151102
$arr[15,1] = 1.459201604194902
$arr[15,2] = 0.1832053746843911
$arr[15,3] = 0.2116479075537185
$arr[15,4] = 1.638171019977889
$arr[15,5] = 1.225857657380459
$arr[15,6] = 0
$arr[15,7] = 0.01542106610589844
$arr[15,8] = 0.7033661723104672
$arr[15,9] = 0.4616512098853889
$arr[15,10] = 0.283300535634794
$arr[15,11] = 0
$arr[15,12] = 0
$arr[15,13] = 0
$arr[15,14] = 0.6210317763557001
$arr[15,15] = 0
$arr[16,0] = 4.023005145363641
$arr[16,1] = 1.419581958746107
$arr[16,2] = 0.1788233247772411
$arr[16,3] = 0.2066894274704296
$arr[16,4] = 1.608486662220187
$arr[16,5] = 1.202509195624515
$arr[16,6] = 0
$arr[16,7] = 0.01409069184156753
$arr[16,8] = 0.6941191711474062
$arr[16,9] = 0.4579804452107012
$arr[16,10] = 0.2767240822117429
$arr[16,11] = 0
$arr[16,12] = 0
$arr[16,13] = 0
$arr[16,14] = 0.6263228672454275
$arr[16,15] = 0
$arr[17,0] = 3.985582376798448
$arr[17,1] = 1.407966861784701
$arr[17,2] = 0.1776043434837362
$arr[17,3] = 0.2050952628878946
$arr[17,4] = 1.596296062330396
$arr[17,5] = 1.192445932509173
$arr[17,6] = 0
$arr[17,7] = 0.01380164613815182
$arr[17,8] = 0.6899083100352073
$arr[17,9] = 0.4554790078535049
$arr[17,10] = 0.2744467516939579
$arr[17,11] = 0
$arr[17,12] = 0
$arr[17,13] = 0
$arr[17,14] = 0.6287692425713871
$arr[17,15] = 0
$arr[18,0] = 4.14275114482615
$arr[18,1] = 1.465943809077999
$arr[18,2] = 0.1839291702955421
$arr[18,3] = 0.2125413020019025
$arr[18,4] = 1.644440373707383
$arr[18,5] = 1.230952023161251
$arr[18,6] = 0
$arr[18,7] = 0.0156201415030699
$arr[18,8] = 0.7054655365676581
$arr[18,9] = 0.462774370549937
$arr[18,10] = 0.2845420115630901
$arr[18,11] = 0
$arr[18,12] = 0
$arr[18,13] = 0
$arr[18,14] = 0.6198431031332063
$arr[18,15] = 0
$arr[19,0] = 4.676917875892627
$arr[19,1] = 1.660933927313522
$arr[19,2] = 0.205168179953688
$arr[19,3] = 0.2378848439292938
$arr[19,4] = 1.812737845783246
$arr[19,5] = 1.366125426628031
$arr[19,6] = 0
$arr[19,7] = 0.02236544638541815
$arr[19,8] = 0.7607216885308787
$arr[19,9] = 0.4900169692829692
$arr[19,10] = 0.3191633309307491
$arr[19,11] = 0
$arr[19,12] = 0
$arr[19,13] = 0
$arr[19,14] = 0.5904430285968658
$arr[19,15] = 0
$arr[20,0] = 5.034448897494485
$arr[20,1] = 1.787017725673365
$arr[20,2] = 0.2188079630610531
$arr[20,3] = 0.2546495540943212
$arr[20,4] = 1.930084427680342
$arr[20,5] = 1.46123352364441
$arr[20,6] = 0
$arr[20,7] = 0.02721630792945451
$arr[20,8] = 0.8002683022411929
$arr[20,9] = 0.5111981085105199
$arr[20,10] = 0.3424289780703447
$arr[20,11] = 0
$arr[20,12] = 0
$arr[20,13] = 0
$arr[20,14] = 0.571050903161364
$arr[20,15] = 0
$arr[21,0] = 4.849792196425938
$arr[21,1] = 1.716432931216616
$arr[21,2] = 0.2110488122217618
$arr[21,3] = 0.2455304126866693
$arr[21,4] = 1.870988268367867
$arr[21,5] = 1.41407985286844
$arr[21,6] = 0
$arr[21,7] = 0.02435027778335197
$arr[21,8] = 0.7809617581987709
$arr[21,9] = 0.5020329151927285
$arr[21,10] = 0.3300651637981389
$arr[21,11] = 0
$arr[21,12] = 0
$arr[21,13] = 0
$arr[21,14] = 0.5801501450706965
$arr[21,15] = 0
$arr[22,0] = 4.143691453699375
$arr[22,1] = 1.458038832879026
$arr[22,2] = 0.1828865809242757
$arr[22,3] = 0.2119204592347543
$arr[22,4] = 1.647632392310825
$arr[22,5] = 1.234697725190955
$arr[22,6] = 0
$arr[22,7] = 0.01515175387634571
$arr[22,8] = 0.7075463662609422
$arr[22,9] = 0.4657722523947001
$arr[22,10] = 0.2841433411880558
$arr[22,11] = 0
$arr[22,12] = 0
$arr[22,13] = 0
$arr[22,14] = 0.6186298952462366
$arr[22,15] = 0
$arr[23,0] = 3.387999071028389
$arr[23,1] = 1.185893868916253
$arr[23,2] = 0.1533731817145707
$arr[23,3] = 0.1769277026405263
$arr[23,4] = 1.419061519118117
$arr[23,5] = 1.051252969036113
$arr[23,6] = 0
$arr[23,7] = 0.00768222593907808
$arr[23,8] = 0.6336719034160865
$arr[23,9] = 0.4306055125643837
$arr[23,10] = 0.2365355343375626
$arr[23,11] = 0
$arr[23,12] = 0
$arr[23,13] = 0
$arr[23,14] = 0.664226785053323
$arr[23,15] = 0

$ws.Range("B2:Q25").Value = $arr